# Week 16 update: a new player (R.Bonnafon) gets logged into both the
# "Rushing" and "Receiving" yards tables. In both sheets this shows up as a
# new column inserted right before the existing "R.Anderson" column (column
# I), shifting R.Anderson and everyone after him one column to the right and
# extending the table from column S to column T.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Insert a new blank column at I, pushing I:S -> J:T.
    $ws.Range("I1").EntireColumn.Insert()

    # Populate the new column's header (row 1) and data row (row 2) the same
    # way every other player column is populated.
    $ws.Range("I1").Value = "R.Bonnafon"
    $ws.Range("I2").Value = "n"
}
